$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7792.0513
$ws.Range("I62").Value = 6035.3335
$ws.Range("J62").Value = 10602.8
$ws.Range("K62").Value = 6035.3335
$ws.Range("L62").Value = 10602.8
$ws.Range("M62").Value = -5411.3335
$ws.Range("N62").Value = -11850.8
$ws.Range("H65").Value = 7792.0513
$ws.Range("I65").Value = 6035.3335
$ws.Range("J65").Value = 10602.8
$ws.Range("K65").Value = 30176.6675
$ws.Range("L65").Value = 53014
$ws.Range("M65").Value = -27056.6675
$ws.Range("N65").Value = -59254
$ws.Range("H76").Value = 3706777.2
$ws.Range("I76").Value = 4632496.5
$ws.Range("J76").Value = 3900.6667
$ws.Range("K76").Value = 4632496.5
$ws.Range("L76").Value = 3900.6667
$ws.Range("M76").Value = -4632181.5
$ws.Range("N76").Value = -4530.6667
$ws.Range("H79").Value = 3706777.2
$ws.Range("I79").Value = 4632496.5
$ws.Range("J79").Value = 3900.6667
$ws.Range("K79").Value = 4632496.5
$ws.Range("L79").Value = 3900.6667
$ws.Range("M79").Value = -4631404.5
$ws.Range("N79").Value = -6084.6667
$ws.Range("H113").Value = 64746.625
$ws.Range("I113").Value = 79087.69500000001
$ws.Range("J113").Value = 2602
$ws.Range("K113").Value = 79087.69500000001
$ws.Range("L113").Value = 2602
$ws.Range("M113").Value = -75833.69500000001
$ws.Range("N113").Value = -9110
$ws.Range("H125").Value = 10192984
$ws.Range("I125").Value = 656
$ws.Range("J125").Value = 18686592
$ws.Range("K125").Value = 5904
$ws.Range("L125").Value = 168179328
$ws.Range("M125").Value = -3444
$ws.Range("N125").Value = -168184248
$ws.Range("H129").Value = 1228.6207
$ws.Range("I129").Value = 471.33334
$ws.Range("J129").Value = 1426.174
$ws.Range("K129").Value = 1414.00002
$ws.Range("L129").Value = 4278.522
$ws.Range("M129").Value = 3585.99998
$ws.Range("N129").Value = -14278.522
$ws.Range("H132").Value = 324830.8
$ws.Range("I132").Value = 369320.9
$ws.Range("J132").Value = 80135.336
$ws.Range("K132").Value = 1107962.7
$ws.Range("L132").Value = 240406.008
$ws.Range("M132").Value = -1105432.7
$ws.Range("N132").Value = -245466.008
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 52928
$ws.Range("I5").Value = 83510.914
$ws.Range("J5").Value = 500.14285
$ws.Range("K5").Value = 83510.914
$ws.Range("L5").Value = 500.14285
$ws.Range("M5").Value = -83398.914
$ws.Range("N5").Value = -724.14285
$ws.Range("H122").Value = 1893.5625
$ws.Range("I122").Value = 1891.9166
$ws.Range("J122").Value = 1898.5
$ws.Range("K122").Value = 5675.7498
$ws.Range("L122").Value = 5695.5
$ws.Range("M122").Value = -3225.7498
$ws.Range("N122").Value = -10595.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 52928
$ws.Range("I4").Value = 83510.914
$ws.Range("J4").Value = 500.14285
$ws.Range("K4").Value = 83510.914
$ws.Range("L4").Value = 500.14285
$ws.Range("M4").Value = -83395.914
$ws.Range("N4").Value = -730.14285
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1542.32
$ws.Range("I31").Value = 876.3077
$ws.Range("J31").Value = 2263.8333
$ws.Range("K31").Value = 876.3077
$ws.Range("L31").Value = 2263.8333
$ws.Range("M31").Value = -581.3077
$ws.Range("N31").Value = -2853.8333
$ws.Range("H34").Value = 1542.32
$ws.Range("I34").Value = 876.3077
$ws.Range("J34").Value = 2263.8333
$ws.Range("K34").Value = 876.3077
$ws.Range("L34").Value = 2263.8333
$ws.Range("M34").Value = -674.3077
$ws.Range("N34").Value = -2667.8333
$ws.Range("H58").Value = 1647.9412
$ws.Range("I58").Value = 638.1081
$ws.Range("J58").Value = 4316.7856
$ws.Range("K58").Value = 638.1081
$ws.Range("L58").Value = 4316.7856
$ws.Range("M58").Value = -435.1081
$ws.Range("N58").Value = -4722.7856
$ws.Range("H99").Value = 20835800
$ws.Range("I99").Value = 62500000
$ws.Range("J99").Value = 3700
$ws.Range("K99").Value = 62500000
$ws.Range("L99").Value = 3700
$ws.Range("M99").Value = -62498502
$ws.Range("N99").Value = -6696
$ws.Range("H126").Value = 20835800
$ws.Range("I126").Value = 62500000
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 187500000
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -187497530
$ws.Range("N126").Value = -16040
$ws.Range("H133").Value = 16607.2
$ws.Range("J133").Value = 16197.333
$ws.Range("L133").Value = 16197.333
$ws.Range("N133").Value = -21257.333
$ws.Range("H136").Value = 1647.9412
$ws.Range("I136").Value = 638.1081
$ws.Range("J136").Value = 4316.7856
$ws.Range("K136").Value = 1914.3243
$ws.Range("L136").Value = 12950.3568
$ws.Range("M136").Value = 635.6756999999998
$ws.Range("N136").Value = -18050.3568
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1747.6578
$ws.Range("I131").Value = 524
$ws.Range("J131").Value = 1933.0605
$ws.Range("K131").Value = 1572
$ws.Range("L131").Value = 5799.181500000001
$ws.Range("M131").Value = 3468
$ws.Range("N131").Value = -15879.1815
$ws.Range("H132").Value = 20834148
$ws.Range("I132").Value = 660
$ws.Range("J132").Value = 55556624
$ws.Range("K132").Value = 5940
$ws.Range("L132").Value = 500009616
$ws.Range("M132").Value = -3410
$ws.Range("N132").Value = -500014676
$ws.Range("H139").Value = 2319.3125
$ws.Range("I139").Value = 2107.7856
$ws.Range("J139").Value = 3800
$ws.Range("K139").Value = 6323.3568
$ws.Range("L139").Value = 11400
$ws.Range("M139").Value = -1183.3568
$ws.Range("N139").Value = -21680
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6243.067
$ws.Range("I70").Value = 6328.6665
$ws.Range("J70").Value = 6043.3335
$ws.Range("K70").Value = 6328.6665
$ws.Range("L70").Value = 6043.3335
$ws.Range("M70").Value = -6058.6665
$ws.Range("N70").Value = -6583.3335
$ws.Range("H73").Value = 6243.067
$ws.Range("I73").Value = 6328.6665
$ws.Range("J73").Value = 6043.3335
$ws.Range("K73").Value = 6328.6665
$ws.Range("L73").Value = 6043.3335
$ws.Range("M73").Value = -5392.6665
$ws.Range("N73").Value = -7915.3335
$ws.Range("H102").Value = 4891.1
$ws.Range("I102").Value = 4752.75
$ws.Range("J102").Value = 4983.3335
$ws.Range("K102").Value = 4752.75
$ws.Range("L102").Value = 4983.3335
$ws.Range("M102").Value = -3130.75
$ws.Range("N102").Value = -8227.333500000001
$ws.Range("H122").Value = 795208.3
$ws.Range("I122").Value = 1588559.4
$ws.Range("J122").Value = 1857.1428
$ws.Range("K122").Value = 4765678.199999999
$ws.Range("L122").Value = 5571.428400000001
$ws.Range("M122").Value = -4763228.199999999
$ws.Range("N122").Value = -10471.4284
$ws.Range("H138").Value = 62999
$ws.Range("J138").Value = 62999
$ws.Range("L138").Value = 62999
$ws.Range("N138").Value = -73279
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3344.889
$ws.Range("I7").Value = 2902.6667
$ws.Range("J7").Value = 3433.3333
$ws.Range("K7").Value = 2902.6667
$ws.Range("L7").Value = 3433.3333
$ws.Range("M7").Value = -2790.6667
$ws.Range("N7").Value = -3657.3333
$ws.Range("H40").Value = 3596.4666
$ws.Range("I40").Value = 2528.1428
$ws.Range("J40").Value = 4531.25
$ws.Range("K40").Value = 2528.1428
$ws.Range("L40").Value = 4531.25
$ws.Range("M40").Value = -2392.1428
$ws.Range("N40").Value = -4803.25
$ws.Range("H122").Value = 3391.0435
$ws.Range("I122").Value = 1978.8
$ws.Range("J122").Value = 3783.3333
$ws.Range("K122").Value = 5936.4
$ws.Range("L122").Value = 11349.9999
$ws.Range("M122").Value = -3486.4
$ws.Range("N122").Value = -16249.9999
$ws.Range("H126").Value = 3344.889
$ws.Range("I126").Value = 2902.6667
$ws.Range("J126").Value = 3433.3333
$ws.Range("K126").Value = 8708.000100000001
$ws.Range("L126").Value = 10299.9999
$ws.Range("M126").Value = -6238.000100000001
$ws.Range("N126").Value = -15239.9999
$ws.Range("H132").Value = 3899.4119
$ws.Range("I132").Value = 2905.5715
$ws.Range("J132").Value = 5504.846
$ws.Range("K132").Value = 8716.7145
$ws.Range("L132").Value = 16514.538
$ws.Range("M132").Value = -6186.7145
$ws.Range("N132").Value = -21574.538
$ws.Range("H133").Value = 52178.25
$ws.Range("J133").Value = 52178.25
$ws.Range("L133").Value = 52178.25
$ws.Range("N133").Value = -57238.25
$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 7000
$ws.Range("K134").Value = 7000
$ws.Range("M134").Value = -1930
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws.Range("H136").Value = 3818.1707
$ws.Range("I136").Value = 2217.2285
$ws.Range("J136").Value = 13157
$ws.Range("K136").Value = 6651.685500000001
$ws.Range("L136").Value = 39471
$ws.Range("M136").Value = -4101.685500000001
$ws.Range("N136").Value = -44571
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0
$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 257200
$ws.Range("J46").Value = 257200
$ws.Range("L46").Value = 257200
$ws.Range("N46").Value = -257662
$ws.Range("H126").Value = 78131.08
$ws.Range("I126").Value = 333767
$ws.Range("J126").Value = 1440.3
$ws.Range("K126").Value = 1001301
$ws.Range("L126").Value = 4320.9
$ws.Range("M126").Value = -998831
$ws.Range("N126").Value = -9260.9
$ws.Range("H134").Value = 257200
$ws.Range("J134").Value = 257200
$ws.Range("L134").Value = 771600
$ws.Range("N134").Value = -776670
$ws.Range("H136").Value = 9288004
$ws.Range("I136").Value = 9834246
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 29502738
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -29500188
$ws.Range("N136").Value = -10800
